$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2430769230769231
$ws.Range("C2").Value = 0.48
$ws.Range("J2").Value = 0.02461538461538462
$ws.Range("P2").Value = 0.1384615384615385
$ws.Range("S2").Value = 0.1138461538461538
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.04294478527607362
$ws.Range("J3").Value = 0.01840490797546012
$ws.Range("P3").Value = 0.754601226993865
$ws.Range("S3").Value = 0.1779141104294479
$ws.Range("J4").Value = 0.05084745762711865
$ws.Range("P4").Value = 0.576271186440678
$ws.Range("S4").Value = 0.3728813559322034
$ws.Range("O5").Value = 0.2
$ws.Range("P5").Value = 0.8
$ws.Range("B6").Value = 0.1005586592178771
$ws.Range("D6").Value = 0.0335195530726257
$ws.Range("F6").Value = 0.0670391061452514
$ws.Range("J6").Value = 0.2625698324022346
$ws.Range("O6").Value = 0.0335195530726257
$ws.Range("Q6").Value = 0.1005586592178771
$ws.Range("R6").Value = 0.0782122905027933
$ws.Range("S6").Value = 0.3240223463687151
$ws.Range("B7").Value = 0.1081081081081081
$ws.Range("D7").Value = 0.01621621621621622
$ws.Range("E7").Value = 0.005405405405405406
$ws.Range("F7").Value = 0.02702702702702703
$ws.Range("J7").Value = 0.1891891891891892
$ws.Range("O7").Value = 0.01081081081081081
$ws.Range("Q7").Value = 0.1621621621621622
$ws.Range("R7").Value = 0.05945945945945946
$ws.Range("S7").Value = 0.4216216216216216
$ws.Range("B8").Value = 0.09550561797752809
$ws.Range("D8").Value = 0.02528089887640449
$ws.Range("F8").Value = 0.04775280898876404
$ws.Range("J8").Value = 0.1573033707865168
$ws.Range("Q8").Value = 0.1769662921348314
$ws.Range("R8").Value = 0.1095505617977528
$ws.Range("S8").Value = 0.3764044943820224
$ws.Range("B9").Value = 0.09146341463414634
$ws.Range("D9").Value = 0.02439024390243903
$ws.Range("E9").Value = 0.006097560975609756
$ws.Range("F9").Value = 0.06707317073170732
$ws.Range("J9").Value = 0.1402439024390244
$ws.Range("O9").Value = 0.02439024390243903
$ws.Range("Q9").Value = 0.1707317073170732
$ws.Range("R9").Value = 0.1158536585365854
$ws.Range("S9").Value = 0.3597560975609756
$ws.Range("B10").Value = 0.1257861635220126
$ws.Range("D10").Value = 0.03066037735849057
$ws.Range("E10").Value = 0.002358490566037736
$ws.Range("F10").Value = 0.06210691823899371
$ws.Range("J10").Value = 0.130503144654088
$ws.Range("O10").Value = 0.01415094339622642
$ws.Range("Q10").Value = 0.2106918238993711
$ws.Range("R10").Value = 0.06918238993710692
$ws.Range("S10").Value = 0.354559748427673
$ws.Range("G11").Value = 0.1374045801526718
$ws.Range("J11").Value = 0.08396946564885496
$ws.Range("K11").Value = 0.1984732824427481
$ws.Range("L11").Value = 0.5725190839694656
$ws.Range("S11").Value = 0.007633587786259542
$ws.Range("G12").Value = 0.7884615384615384
$ws.Range("J12").Value = 0.1730769230769231
$ws.Range("K12").Value = 0.00641025641025641
$ws.Range("L12").Value = 0.02564102564102564
$ws.Range("S12").Value = 0.00641025641025641
$ws.Range("G13").Value = 0.6326530612244898
$ws.Range("J13").Value = 0.3265306122448979
$ws.Range("S13").Value = 0.04081632653061224
$ws.Range("F15").Value = 0.01530612244897959
$ws.Range("H15").Value = 0.1479591836734694
$ws.Range("I15").Value = 0.05612244897959184
$ws.Range("J15").Value = 0.4183673469387755
$ws.Range("K15").Value = 0.07653061224489796
$ws.Range("M15").Value = 0.01530612244897959
$ws.Range("O15").Value = 0.03571428571428571
$ws.Range("S15").Value = 0.2346938775510204
$ws.Range("F16").Value = 0.01477832512315271
$ws.Range("H16").Value = 0.1477832512315271
$ws.Range("I16").Value = 0.07881773399014778
$ws.Range("J16").Value = 0.4630541871921182
$ws.Range("K16").Value = 0.1083743842364532
$ws.Range("M16").Value = 0.02463054187192118
$ws.Range("O16").Value = 0.06403940886699508
$ws.Range("S16").Value = 0.09852216748768473
$ws.Range("F17").Value = 0.01732673267326733
$ws.Range("H17").Value = 0.1658415841584159
$ws.Range("I17").Value = 0.08663366336633663
$ws.Range("J17").Value = 0.4554455445544555
$ws.Range("K17").Value = 0.07673267326732673
$ws.Range("M17").Value = 0.02722772277227723
$ws.Range("O17").Value = 0.05198019801980198
$ws.Range("S17").Value = 0.1188118811881188
$ws.Range("F18").Value = 0.01775147928994083
$ws.Range("H18").Value = 0.1301775147928994
$ws.Range("I18").Value = 0.07100591715976332
$ws.Range("J18").Value = 0.514792899408284
$ws.Range("K18").Value = 0.09467455621301775
$ws.Range("M18").Value = 0.01775147928994083
$ws.Range("O18").Value = 0.05325443786982249
$ws.Range("S18").Value = 0.1005917159763314
$ws.Range("F19").Value = 0.01515151515151515
$ws.Range("H19").Value = 0.1871657754010695
$ws.Range("I19").Value = 0.08288770053475936
$ws.Range("J19").Value = 0.3894830659536542
$ws.Range("K19").Value = 0.1051693404634581
$ws.Range("M19").Value = 0.02584670231729055
$ws.Range("O19").Value = 0.0748663101604278
$ws.Range("S19").Value = 0.1194295900178253
